$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update file paths to be relative (engine='openpyxl' change in daily_map.py)
$ws.Range("B4").Value = "data/dtm_depth_padded.tif"
$ws.Range("B5").Value = "data/depth_extended.tif"
$ws.Range("B6").Value = "data/199_canalblocks_20191008b.shp"
$ws.Range("B7").Value = "data/weather_station_coordinates.xlsx"
$ws.Range("B8").Value = "data/new_area/mesh_0.02.msh"
$ws.Range("B9").Value = "data/new_area/canal_network_matrix_50meters.p"
$ws.Range("B10").Value = "data/dtm_big_area_depth_padded.tif"

# Update the active selection to D16
$ws.Range("D16").Select()
